# Updated cryptos list (automated price/volume refresh).
# Column D ("Price") cells are stored as plain text in this workbook (e.g.
# "235.51", "1.832.70"), matching the original inlineStr cell type. Values
# that look like a bare number are prefixed with a leading apostrophe so
# Excel's COM layer stores them as text instead of silently coercing them
# into Double values (which would corrupt values like "0.05500" or
# "1.147.35" and introduce floating-point noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.299.48"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.830.48"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "'235.51"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").Value = "'0.6028"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'0.06979"
$ws.Range("E8").Value = "  -4.94%  "
$ws.Range("D9").Value = "'0.2771"
$ws.Range("E9").Value = "  -3.72%  "
$ws.Range("D10").Value = "'23.62"
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("D11").Value = "'0.07613"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "1.843.32"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'4.758"
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").Value = "'0.6328"
$ws.Range("E14").Value = "  -4.12%  "
$ws.Range("D15").Value = "'0.000009861"
$ws.Range("E15").Value = "  -4.41%  "
$ws.Range("D16").Value = "'77.88"
$ws.Range("E16").Value = "  -4.29%  "
$ws.Range("D17").Value = "29.013.76"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "'5.598"
$ws.Range("E18").Value = "  -10.24%  "
$ws.Range("D19").Value = "'218.02"
$ws.Range("E19").Value = "  -7.87%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "'11.61"
$ws.Range("E21").Value = "  -4.77%  "
$ws.Range("D22").Value = "'6.918"
$ws.Range("E22").Value = "  -4.03%  "
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Value = "'156.58"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'7.988"
$ws.Range("E25").Value = "  -5.02%  "
$ws.Range("D26").Value = "'0.1294"
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("D27").Value = "'16.56"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").Value = "'0.06432"
$ws.Range("E28").Value = "  -7.05%  "
$ws.Range("E29").Value = "  -3.27%  "
$ws.Range("D30").Value = "'1.445"
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("D31").Value = "'3.841"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").Value = "'3.805"
$ws.Range("E32").Value = "  -5.19%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.097"
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'1.731"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").Value = "'0.6494"
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").Value = "'2.546"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").Value = "'2.758"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").Value = "'0.01758"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").Value = "'6.601"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").Value = "1.144.70"
$ws.Range("E40").Value = "  -6.98%  "
$ws.Range("D41").Value = "'0.8950"
$ws.Range("E41").Value = "  -5.29%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "2.001.74"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'100.96"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "'62.35"
$ws.Range("E45").Value = "  -4.22%  "
$ws.Range("D46").Value = "'0.00000000114"
$ws.Range("E46").Value = "  -3.49%  "
$ws.Range("D47").Value = "'1.623"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("D48").Value = "'8.501"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4547"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05500"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").Value = "'6.402"
$ws.Range("E51").Value = "  -6.94%  "
